$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 3
$ws.Range("H2").Value = 3.3
$ws.Range("I2").Value = 2.45
$ws.Range("K2").Value = 2.05
$ws.Range("L2").Value = 3.2
$ws.Range("S2").Value = 1.44
$ws.Range("T2").Value = 2.63
$ws.Range("W2").Value = 8.5
$ws.Range("AB2").Value = 34
$ws.Range("AC2").Value = 8.5
$ws.Range("AE2").Value = 15
$ws.Range("AG2").Value = 351
$ws.Range("AH2").Value = 7.5
$ws.Range("AS2").Value = 201
$ws.Range("AT2").Value = 2.63
$ws.Range("AU2").Value = 8
$ws.Range("AV2").Value = 51
$ws.Range("BA2").Value = 67
$ws.Range("BC2").Value = 151

# Row 4
$ws.Range("G4").Value = 2.88
$ws.Range("J4").Value = 3.6
$ws.Range("L4").Value = 3.5
$ws.Range("Q4").Value = 2.5
$ws.Range("R4").Value = 1.5
$ws.Range("U4").Value = 2
$ws.Range("V4").Value = 1.75
$ws.Range("X4").Value = 13
$ws.Range("Y4").Value = 12
$ws.Range("AC4").Value = 6
$ws.Range("AD4").Value = 5.5
$ws.Range("AG4").Value = 451

# Row 5
$ws.Range("G5").Value = 1.85
$ws.Range("H5").Value = 3.2
$ws.Range("I5").Value = 5
$ws.Range("J5").Value = 2.63
$ws.Range("Q5").Value = 2.88
$ws.Range("R5").Value = 1.4
$ws.Range("X5").Value = 7
$ws.Range("AW5").Value = 6.5

# Row 7
$ws.Range("M7").Value = 1.08
$ws.Range("N7").Value = 8
$ws.Range("O7").Value = 1.4
$ws.Range("P7").Value = 2.75
$ws.Range("Q7").Value = 2.25
$ws.Range("R7").Value = 1.62

# Row 8
$ws.Range("G8").Value = 2.3
$ws.Range("I8").Value = 3.4
$ws.Range("J8").Value = 3.1
$ws.Range("M8").Value = 1.11
$ws.Range("N8").Value = 6.5
$ws.Range("X8").Value = 9.5
$ws.Range("Z8").Value = 21
$ws.Range("AA8").Value = 21
$ws.Range("AH8").Value = 8
$ws.Range("AK8").Value = 41
$ws.Range("AN8").Value = 4
$ws.Range("AO8").Value = 13
$ws.Range("AX8").Value = 21

# Row 9
$ws.Range("G9").Value = 2.05
$ws.Range("H9").Value = 3
$ws.Range("I9").Value = 4.2
$ws.Range("J9").Value = 2.88
$ws.Range("L9").Value = 4.75
$ws.Range("O9").Value = 1.5
$ws.Range("P9").Value = 2.4
$ws.Range("Q9").Value = 2.6
$ws.Range("R9").Value = 1.48
$ws.Range("X9").Value = 8
$ws.Range("Y9").Value = 9.5
$ws.Range("Z9").Value = 17
$ws.Range("AA9").Value = 21
$ws.Range("AD9").Value = 6
$ws.Range("AH9").Value = 9
$ws.Range("AI9").Value = 19
$ws.Range("AN9").Value = 3.75
$ws.Range("AO9").Value = 12
$ws.Range("AQ9").Value = 41
$ws.Range("AS9").Value = 251
$ws.Range("AT9").Value = 2.25
$ws.Range("AX9").Value = 23
$ws.Range("BA9").Value = 151

# Row 10
$ws.Range("M10").Value = 1.1
$ws.Range("N10").Value = 7
$ws.Range("Q10").Value = 2.6
$ws.Range("R10").Value = 1.48
$ws.Range("AL10").Value = 23
$ws.Range("AO10").Value = 19
$ws.Range("AQ10").Value = 67
$ws.Range("AW10").Value = 4.33
$ws.Range("BB10").Value = 251

# Row 11
$ws.Range("H11").Value = 3.15
$ws.Range("I11").Value = 5.1
$ws.Range("J11").Value = 2.32
$ws.Range("K11").Value = 2.05
$ws.Range("L11").Value = 5.3
$ws.Range("M11").Value = 1.09
$ws.Range("N11").Value = 6.1
$ws.Range("O11").Value = 1.4
$ws.Range("P11").Value = 2.72
$ws.Range("Q11").Value = 2.2
$ws.Range("R11").Value = 1.62
$ws.Range("S11").Value = 1.45
$ws.Range("T11").Value = 2.57
$ws.Range("U11").Value = 2
$ws.Range("V11").Value = 1.72
$ws.Range("W11").Value = 5.6
$ws.Range("X11").Value = 7.4
$ws.Range("Y11").Value = 8.25
$ws.Range("AA11").Value = 15.5
$ws.Range("AB11").Value = 32
$ws.Range("AC11").Value = 6.1
$ws.Range("AD11").Value = 6.2
$ws.Range("AE11").Value = 17
$ws.Range("AF11").Value = 100
$ws.Range("AG11").Value = 900
$ws.Range("AH11").Value = 11.5
$ws.Range("AI11").Value = 29
$ws.Range("AJ11").Value = 16.5
$ws.Range("AL11").Value = 60
$ws.Range("AN11").Value = 3.5
$ws.Range("AO11").Value = 8.75
$ws.Range("AP11").Value = 18.5
$ws.Range("AQ11").Value = 32
$ws.Range("AR11").Value = 65
$ws.Range("AT11").Value = 2.57
$ws.Range("AU11").Value = 7.4
$ws.Range("AV11").Value = 70
$ws.Range("AW11").Value = 6.7

# Row 12
$ws.Range("G12").Value = 2.65
$ws.Range("H12").Value = 2.57
$ws.Range("I12").Value = 3.2
$ws.Range("K12").Value = 1.85
$ws.Range("L12").Value = 3.8
$ws.Range("M12").Value = 1.13
$ws.Range("N12").Value = 5.1
$ws.Range("O12").Value = 1.5
$ws.Range("P12").Value = 2.42
$ws.Range("Q12").Value = 2.5
$ws.Range("R12").Value = 1.47
$ws.Range("S12").Value = 1.55
$ws.Range("T12").Value = 2.37
$ws.Range("U12").Value = 1.93
$ws.Range("V12").Value = 1.78
$ws.Range("W12").Value = 6.6
$ws.Range("AA12").Value = 26
$ws.Range("AB12").Value = 40
$ws.Range("AC12").Value = 5.1
$ws.Range("AE12").Value = 14.5
$ws.Range("AF12").Value = 80
$ws.Range("AG12").Value = 800
$ws.Range("AH12").Value = 7.5
$ws.Range("AJ12").Value = 11.25
$ws.Range("AL12").Value = 35
$ws.Range("AM12").Value = 45
$ws.Range("AO12").Value = 15
$ws.Range("AP12").Value = 23
$ws.Range("AR12").Value = 110
$ws.Range("AT12").Value = 2.32
$ws.Range("AU12").Value = 6.7
$ws.Range("AV12").Value = 65
$ws.Range("AX12").Value = 19
$ws.Range("AY12").Value = 26
$ws.Range("AZ12").Value = 100
$ws.Range("BA12").Value = 150
$ws.Range("BB12").Value = 350

# Row 14
$ws.Range("K14").Value = 1.95
$ws.Range("M14").Value = 1.08
$ws.Range("N14").Value = 7.5
$ws.Range("S14").Value = 1.53
$ws.Range("T14").Value = 2.38
$ws.Range("U14").Value = 2.1
$ws.Range("V14").Value = 1.67
$ws.Range("X14").Value = 9.5
$ws.Range("Y14").Value = 10
$ws.Range("AE14").Value = 19
$ws.Range("AG14").Value = 501
$ws.Range("AN14").Value = 4
$ws.Range("AQ14").Value = 51
$ws.Range("AS14").Value = 251
$ws.Range("AT14").Value = 2.38
$ws.Range("AU14").Value = 9
$ws.Range("AX14").Value = 21
$ws.Range("BB14").Value = 301

# Row 15
$ws.Range("G15").Value = 4.33
$ws.Range("H15").Value = 3.5
$ws.Range("I15").Value = 1.8
$ws.Range("K15").Value = 2.05
$ws.Range("L15").Value = 2.5
$ws.Range("M15").Value = 1.08
$ws.Range("N15").Value = 8
$ws.Range("W15").Value = 10
$ws.Range("AC15").Value = 8
$ws.Range("AR15").Value = 151
$ws.Range("AW15").Value = 3.6
$ws.Range("AZ15").Value = 34

# Row 16
$ws.Range("G16").Value = 3
$ws.Range("H16").Value = 3.5
$ws.Range("I16").Value = 2.25
$ws.Range("J16").Value = 3.4
$ws.Range("O16").Value = 1.17
$ws.Range("P16").Value = 5
$ws.Range("Q16").Value = 1.6
$ws.Range("R16").Value = 2.3
$ws.Range("AB16").Value = 26
$ws.Range("AK16").Value = 21
$ws.Range("AW16").Value = 4.5

# Row 18
$ws.Range("M18").Value = 1.05
$ws.Range("N18").Value = 11
$ws.Range("Q18").Value = 1.85
$ws.Range("R18").Value = 2

# Row 19
$ws.Range("G19").Value = 1.4
$ws.Range("H19").Value = 4.75
$ws.Range("I19").Value = 6
$ws.Range("J19").Value = 1.91
$ws.Range("K19").Value = 2.38
$ws.Range("M19").Value = 1.03
$ws.Range("N19").Value = 10
$ws.Range("AA19").Value = 12
$ws.Range("AC19").Value = 13
$ws.Range("AD19").Value = 9.5
$ws.Range("AM19").Value = 41
$ws.Range("AQ19").Value = 19

# Row 21
$ws.Range("G21").Value = 6.25
$ws.Range("H21").Value = 4.5
$ws.Range("J21").Value = 5.5
$ws.Range("L21").Value = 1.91
$ws.Range("AD21").Value = 9
$ws.Range("AG21").Value = 450
$ws.Range("AI21").Value = 8
$ws.Range("AO21").Value = 29
$ws.Range("AR21").Value = 101
$ws.Range("AS21").Value = 350
$ws.Range("AZ21").Value = 19
